# Fruta / hortaliza, semanal
# Insert a new weekly price record row at row 334, shifting the existing
# rows 334:355 down to 335:356.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 334 (pushes old 334..355 -> 335..356)
$ws.Rows.Item(334).Insert()

# Populate the new row 334 with the new record's data
$ws.Range("A334").Value = 9
$ws.Range("B334").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C334").Value = "Metropolitana"
$ws.Range("D334").Value = 45021
$ws.Range("E334").Value = 13
$ws.Range("F334").Value = 100112001
$ws.Range("G334").Value = "Berenjena"
$ws.Range("H334").Value = "Sin especificar"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 65
$ws.Range("K334").Value = 8000
$ws.Range("L334").Value = 8000
$ws.Range("M334").Value = 8000
$ws.Range("N334").Value = "$/caja 50 unidades"
$ws.Range("O334").Value = "Región de Arica y Parinacota"
$ws.Range("P334").Value = 160
$ws.Range("Q334").Value = 50
$ws.Range("R334").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D
$ws.Range("D334").NumberFormat = $ws.Range("D335").NumberFormat
